$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44335
$ws.Range("J2").Value = 150
$ws.Range("M2").Value = 633
$ws.Range("P2").Value = 106

# Row 3
$ws.Range("D3").Value = 44335
$ws.Range("J3").Value = 50

# Row 4
$ws.Range("D4").Value = 44188

# Row 5
$ws.Range("D5").Value = 44188

# Row 8
$ws.Range("D8").Value = 44293

# Row 9
$ws.Range("D9").Value = 44293

# Row 10
$ws.Range("D10").Value = 44321

# Row 11
$ws.Range("D11").Value = 44321

# Row 12
$ws.Range("D12").Value = 44328
$ws.Range("J12").Value = 100

# Row 13
$ws.Range("D13").Value = 44328
$ws.Range("J13").Value = 50

# Row 14
$ws.Range("D14").Value = 44308
$ws.Range("J14").Value = 200

# Row 15
$ws.Range("D15").Value = 44308
$ws.Range("J15").Value = 100

# Row 16
$ws.Range("D16").Value = 44358
$ws.Range("J16").Value = 200
$ws.Range("M16").Value = 650
$ws.Range("P16").Value = 108

# Row 17
$ws.Range("D17").Value = 44358
$ws.Range("J17").Value = 100
